$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.455.85"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "3.052.02"
$ws.Range("E3").Value = "  +2.44%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'385.17"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").Value = "'103.11"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  -0.48%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.585"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("D10").Value = "'36.79"
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D11").Value = "'0.138"
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "'0.0862"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "3.531.06"
$ws.Range("E13").Value = "  +2.47%  "
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").Value = "'7.77"
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").Value = "3.053.29"
$ws.Range("E16").Value = "  +2.06%  "
$ws.Range("D17").Value = "'0.971"
$ws.Range("E17").Value = "  -2.83%  "
$ws.Range("D18").Value = "'10.59"
$ws.Range("E18").Value = "  -5.34%  "
$ws.Range("D19").Value = "51.526.64"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").Value = "'3.16"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("D21").Value = "'12.40"
$ws.Range("E21").Value = "  -1.20%  "
$ws.Range("D22").Value = "0.0₃0966"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("D23").Value = "'70.13"
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").Value = "'268.09"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").Value = "'3.15"
$ws.Range("E25").Value = "  -2.22%  "
$ws.Range("D26").Value = "'8.21"
$ws.Range("E26").Value = "  +4.20%  "
$ws.Range("D27").Value = "'26.92"
$ws.Range("E27").Value = "  +3.25%  "
$ws.Range("E28").Value = "  +3.04%  "
$ws.Range("D29").Value = "'7.24"
$ws.Range("E29").Value = "  -3.90%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("D33").Value = "'34.86"
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "'50.44"
$ws.Range("E35").Value = "  -2.00%  "
$ws.Range("D36").Value = "'0.0447"
$ws.Range("E36").Value = "  +1.76%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  +2.41%  "
$ws.Range("D39").Value = "'0.294"
$ws.Range("E39").Value = "  +8.22%  "
$ws.Range("D40").Value = "'16.97"
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("D41").Value = "'1.87"
$ws.Range("E41").Value = "  +1.18%  "
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").Value = "'2.55"
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").Value = "'124.85"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").Value = "'3.77"
$ws.Range("E45").Value = "  +3.07%  "
$ws.Range("D46").Value = "'21.87"
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("E47").Value = "  +3.02%  "
$ws.Range("D48").Value = "'2.39"
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("D49").Value = "2.030.86"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").Value = "3.348.23"
$ws.Range("E50").Value = "  +2.54%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.205"
$ws.Range("E51").Value = "  +6.03%  "
